# Scheduled runner update: refresh market price / profit columns (H-N) per leve row
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 2266.6667
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 2900
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 2900
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -3038

# Row 62
$ws.Range("H62").Value = 7139.3335
$ws.Range("I62").Value = 6220.5
$ws.Range("J62").Value = 8189.4287
$ws.Range("K62").Value = 6220.5
$ws.Range("L62").Value = 8189.4287
$ws.Range("M62").Value = -5596.5
$ws.Range("N62").Value = -9437.4287

# Row 65
$ws.Range("H65").Value = 7139.3335
$ws.Range("I65").Value = 6220.5
$ws.Range("J65").Value = 8189.4287
$ws.Range("K65").Value = 31102.5
$ws.Range("L65").Value = 40947.14350000001
$ws.Range("M65").Value = -27982.5
$ws.Range("N65").Value = -47187.14350000001

# Row 100
$ws.Range("H100").Value = 4832.091
$ws.Range("I100").Value = 3815.3
$ws.Range("K100").Value = 3815.3
$ws.Range("M100").Value = -3274.3

# Row 103
$ws.Range("H103").Value = 1151.6471
$ws.Range("J103").Value = 768.8889
$ws.Range("L103").Value = 2306.6667
$ws.Range("N103").Value = -3478.6667

# Row 107
$ws.Range("H107").Value = 1389.3158
$ws.Range("I107").Value = 1540.6471
$ws.Range("K107").Value = 1540.6471
$ws.Range("M107").Value = 379.3529000000001

# Row 113
$ws.Range("H113").Value = 6398.222
$ws.Range("I113").Value = 5821.6665
$ws.Range("J113").Value = 6686.5
$ws.Range("K113").Value = 5821.6665
$ws.Range("L113").Value = 6686.5
$ws.Range("M113").Value = -2567.6665
$ws.Range("N113").Value = -13194.5

# Row 125
$ws.Range("H125").Value = 5727737.5
$ws.Range("J125").Value = 1660.625
$ws.Range("L125").Value = 14945.625
$ws.Range("N125").Value = -19865.625

# Row 135
$ws.Range("H135").Value = 3342.2307
$ws.Range("I135").Value = 3257.375
$ws.Range("K135").Value = 29316.375
$ws.Range("M135").Value = -26781.375

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 629.9091
$ws.Range("I2").Value = 410.85715
$ws.Range("K2").Value = 410.85715
$ws.Range("M2").Value = -297.85715

# Row 32
$ws.Range("H32").Value = 21657.453
$ws.Range("I32").Value = 22165.117
$ws.Range("K32").Value = 22165.117
$ws.Range("M32").Value = -21878.117

# Row 61
$ws.Range("H61").Value = 10540.5
$ws.Range("I61").Value = 6735.385
$ws.Range("K61").Value = 6735.385
$ws.Range("M61").Value = -6523.385

# Row 97
$ws.Range("H97").Value = 700.1429000000001
$ws.Range("J97").Value = 929.5
$ws.Range("L97").Value = 929.5
$ws.Range("N97").Value = -1921.5

# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# Row 116
$ws.Range("H116").Value = 629.9091
$ws.Range("I116").Value = 410.85715
$ws.Range("K116").Value = 410.85715
$ws.Range("M116").Value = 1883.14285

# Row 136
$ws.Range("H136").Value = 10540.5
$ws.Range("I136").Value = 6735.385
$ws.Range("K136").Value = 20206.155
$ws.Range("M136").Value = -17656.155

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 629.9091
$ws.Range("I3").Value = 410.85715
$ws.Range("K3").Value = 410.85715
$ws.Range("M3").Value = -296.85715

# Row 80
$ws.Range("H80").Value = 690.93335
$ws.Range("I80").Value = 992
$ws.Range("J80").Value = 346.85715
$ws.Range("K80").Value = 992
$ws.Range("L80").Value = 346.85715
$ws.Range("M80").Value = 6
$ws.Range("N80").Value = -2342.85715

# Row 83
$ws.Range("H83").Value = 690.93335
$ws.Range("I83").Value = 992
$ws.Range("J83").Value = 346.85715
$ws.Range("K83").Value = 4960
$ws.Range("L83").Value = 1734.28575
$ws.Range("M83").Value = 32
$ws.Range("N83").Value = -11718.28575

# Row 86
$ws.Range("H86").Value = 1944.091
$ws.Range("I86").Value = 1738.5
$ws.Range("K86").Value = 1738.5
$ws.Range("M86").Value = -615.5

# Row 89
$ws.Range("H89").Value = 1944.091
$ws.Range("I89").Value = 1738.5
$ws.Range("K89").Value = 8692.5
$ws.Range("M89").Value = -3076.5

# Row 94
$ws.Range("H94").Value = 2778234.8
$ws.Range("I94").Value = 483.44446
$ws.Range("J94").Value = 11111489
$ws.Range("K94").Value = 483.44446
$ws.Range("L94").Value = 11111489
$ws.Range("M94").Value = -32.44445999999999
$ws.Range("N94").Value = -11112391

# Row 99
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 7225.2354
$ws.Range("I99").Value = 7184.5
$ws.Range("J99").Value = 7283.4287
$ws.Range("K99").Value = 7184.5
$ws.Range("L99").Value = 7283.4287
$ws.Range("M99").Value = -5686.5
$ws.Range("N99").Value = -10279.4287

# Row 126
$ws.Range("H126").Value = 7225.2354
$ws.Range("I126").Value = 7184.5
$ws.Range("J126").Value = 7283.4287
$ws.Range("K126").Value = 21553.5
$ws.Range("L126").Value = 21850.2861
$ws.Range("M126").Value = -19083.5
$ws.Range("N126").Value = -26790.2861

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 4334.154
$ws.Range("J68").Value = 4643
$ws.Range("L68").Value = 13929
$ws.Range("N68").Value = -15551

# Row 71
$ws.Range("H71").Value = 4334.154
$ws.Range("J71").Value = 4643
$ws.Range("L71").Value = 41787
$ws.Range("N71").Value = -49899

# Row 125
$ws.Range("H125").Value = 5800
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

# Row 132
$ws.Range("H132").Value = 1503.2106
$ws.Range("I132").Value = 1531.7142
$ws.Range("J132").Value = 1486.5834
$ws.Range("K132").Value = 13785.4278
$ws.Range("L132").Value = 13379.2506
$ws.Range("M132").Value = -11255.4278
$ws.Range("N132").Value = -18439.2506

$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 70000
$ws.Range("J62").Value = 70000
$ws.Range("L62").Value = 70000
$ws.Range("N62").Value = -71372

# Row 65
$ws.Range("H65").Value = 70000
$ws.Range("J65").Value = 70000
$ws.Range("L65").Value = 210000
$ws.Range("N65").Value = -216864

# Row 97
$ws.Range("H97").Value = 725.89655
$ws.Range("I97").Value = 674.93335
$ws.Range("J97").Value = 780.5
$ws.Range("K97").Value = 674.93335
$ws.Range("L97").Value = 780.5
$ws.Range("M97").Value = -178.93335
$ws.Range("N97").Value = -1772.5

# Row 102
$ws.Range("H102").Value = 1808.7
$ws.Range("I102").Value = 1386
$ws.Range("J102").Value = 3499.5
$ws.Range("K102").Value = 1386
$ws.Range("L102").Value = 3499.5
$ws.Range("M102").Value = 236
$ws.Range("N102").Value = -6743.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 4181.0435
$ws.Range("I46").Value = 2249.5
$ws.Range("J46").Value = 4365
$ws.Range("K46").Value = 2249.5
$ws.Range("L46").Value = 4365
$ws.Range("M46").Value = -2061.5
$ws.Range("N46").Value = -4741

# Row 61
$ws.Range("H61").Value = 2084.8333
$ws.Range("I61").Value = 2101.8
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 2101.8
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1899.8
$ws.Range("N61").Value = -2404

# Row 69
$ws.Range("H69").Value = 146666.67
$ws.Range("J69").Value = 146666.67
$ws.Range("L69").Value = 146666.67
$ws.Range("N69").Value = -148288.67

# Row 72
$ws.Range("H72").Value = 146666.67
$ws.Range("J72").Value = 146666.67
$ws.Range("L72").Value = 440000.01
$ws.Range("N72").Value = -448112.01

# Row 107
$ws.Range("H107").Value = 2966.3333
$ws.Range("I107").Value = 2966.3333
$ws.Range("K107").Value = 2966.3333
$ws.Range("M107").Value = -1046.3333

# Row 113
$ws.Range("H113").Value = 2084.8333
$ws.Range("I113").Value = 2101.8
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2101.8
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 68.19999999999982
$ws.Range("N113").Value = -6340

# Row 140
$ws.Range("H140").Value = 105407.25
$ws.Range("J140").Value = 105407.25
$ws.Range("L140").Value = 105407.25
$ws.Range("N140").Value = -115767.25

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4566.838
$ws.Range("I81").Value = 4194.485
$ws.Range("K81").Value = 8388.969999999999
$ws.Range("M81").Value = -7327.969999999999

# Row 84
$ws.Range("H84").Value = 4566.838
$ws.Range("I84").Value = 4194.485
$ws.Range("K84").Value = 41944.85
$ws.Range("M84").Value = -36640.85

# Row 113
$ws.Range("H113").Value = 572.2045000000001
$ws.Range("I113").Value = 642.0625
$ws.Range("J113").Value = 385.91666
$ws.Range("K113").Value = 1926.1875
$ws.Range("L113").Value = 1157.74998
$ws.Range("M113").Value = 243.8125
$ws.Range("N113").Value = -5497.749980000001

# Row 122
$ws.Range("H122").Value = 4032.9546
$ws.Range("I122").Value = 4240.727
$ws.Range("K122").Value = 12722.181
$ws.Range("M122").Value = -10272.181
